# ----------------------------------------------------------------------------
# Wealth by Race over Lifetime - "updated text, final edits, ready for QC"
# ----------------------------------------------------------------------------
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Worksheet (tab) name: capitalization fix "By" -> "by"
$ws.Name = "Wealth by Race over Lifetime"

# 2. Chart-1 header: hyphen -> en dash in the birth-year range
$avgHeader = $ws.Range("A1")
$avgHeader.Value = "Average Family Wealth for Those Born 1943" + [char]0x2013 + "51 by Race"

# 3. Notes text: reword the trailing sentence
#    "...age 65-73. In 1983, ages 32-40."
#      -> "...age 65-73; in 1983, they were ages 32-40."  (with en dashes)
#    Keep the bold "Notes:" lead-in and the regular formatting of the body
#    text intact (Characters() rewrites reset formatting, so both runs are
#    re-stamped explicitly afterwards).
$notes = $ws.Range("A12")

$notesBody = $notes.Characters(8, 9999)
$notesBody.Text = "2016 dollars. Hispanic sample size too small to show. Age is defined as the age of the household head. In 2016, these people were age 65" + [char]0x2013 + "73; in 1983, they were ages 32" + [char]0x2013 + "40."

$notesBodyFmt = $notes.Characters(7, 9999)
$notesBodyFmt.Font.Name = "Calibri"
$notesBodyFmt.Font.Size = 11
$notesBodyFmt.Font.Bold = $false
$notesBodyFmt.Font.Color = 0

$notesPrefix = $notes.Characters(1, 6)
$notesPrefix.Font.Name = "Calibri"
$notesPrefix.Font.Size = 11
$notesPrefix.Font.Bold = $true
$notesPrefix.Font.Color = 0

# 4. Leave the selection where the author left it when the file was saved for QC
$null = $ws.Range("C21").Select()
